$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$plainRange = $ws.Range("A2,C2:E2,G2:H2")
foreach ($area in $plainRange.Areas) {
    $area.Font.Name = "Calibri"
    $area.Font.Size = 10
    $area.Font.ThemeColor = 1
}

$ws.Range("F2").Font.Name = "Calibri"
$ws.Range("F2").Font.Size = 10
$ws.Range("F2").Font.ThemeColor = 1

$ws.Range("A2").Value = "MCH223-1"
$ws.Range("C2").Value = "CALENDAR"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24F | GRAP COUNT NUMER: NONE"
